$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (CC1:CV1) for the 20 additional feature-importance columns ---
  $ws.Cells.Item(1, 81).Value = "diameters importance"
  $ws.Cells.Item(1, 82).Value = "ratio_area_over_diameter importance"
  $ws.Cells.Item(1, 83).Value = "ratio_perimeter_over_diameter importance"
  $ws.Cells.Item(1, 84).Value = "Dense Urban importance"
  $ws.Cells.Item(1, 85).Value = "Industrial importance"
  $ws.Cells.Item(1, 86).Value = "N,A importance"
  $ws.Cells.Item(1, 87).Value = "Rural importance"
  $ws.Cells.Item(1, 88).Value = "Sparse Urban importance"
  $ws.Cells.Item(1, 89).Value = "Urban Slum importance"
  $ws.Cells.Item(1, 90).Value = "Barren Land importance"
  $ws.Cells.Item(1, 91).Value = "Coastal importance"
  $ws.Cells.Item(1, 92).Value = "Dense Forest importance"
  $ws.Cells.Item(1, 93).Value = "Desert importance"
  $ws.Cells.Item(1, 94).Value = "Farms importance"
  $ws.Cells.Item(1, 95).Value = "Grass Land importance"
  $ws.Cells.Item(1, 96).Value = "Hills importance"
  $ws.Cells.Item(1, 97).Value = "Lakes importance"
  $ws.Cells.Item(1, 98).Value = "River importance"
  $ws.Cells.Item(1, 99).Value = "Snow importance"
  $ws.Cells.Item(1, 100).Value = "Sparse Forest importance"

# Copy the header style (bold, centered, bordered) from the last existing header cell
$ws.Range("CB1").Copy()
$ws.Range("CC1:CV1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New data rows (9, 10, 11) ---
  $ws.Cells.Item(9, 1).Value = 0.9999944350044395
  $ws.Cells.Item(9, 2).Value = 10
  $ws.Cells.Item(9, 3).Value = 0.5501898162519842
  $ws.Cells.Item(9, 4).Value = 0.01069409847699901
  $ws.Cells.Item(9, 5).Value = 0.5643725361706374
  $ws.Cells.Item(9, 6).Value = 0.5277711643627568
  $ws.Cells.Item(9, 7).Value = 0.04046192694419542
  $ws.Cells.Item(9, 8).Value = 0.001464931012790026
  $ws.Cells.Item(9, 9).Value = 0.02663389284411285
  $ws.Cells.Item(9, 65).Value = 0.03578403188927403
  $ws.Cells.Item(9, 66).Value = 0.04770310805630882
  $ws.Cells.Item(9, 67).Value = 0.06887897750582037
  $ws.Cells.Item(9, 68).Value = 0.06438403119721645
  $ws.Cells.Item(9, 69).Value = 0.06463231681652754
  $ws.Cells.Item(9, 70).Value = 0.03078707204472345
  $ws.Cells.Item(9, 71).Value = 0.02704967297420518
  $ws.Cells.Item(9, 72).Value = 0.04392286493728839
  $ws.Cells.Item(9, 73).Value = 0.02828644428033223
  $ws.Cells.Item(9, 74).Value = 0.02998898036313816
  $ws.Cells.Item(9, 75).Value = 0.035634589507616
  $ws.Cells.Item(9, 76).Value = 0.05346282437410078
  $ws.Cells.Item(9, 77).Value = 0.06300635678054092
  $ws.Cells.Item(9, 78).Value = 0.04510933351600396
  $ws.Cells.Item(9, 79).Value = 0.0956296260451878
  $ws.Cells.Item(9, 80).Value = 0.05559398111673422
  $ws.Cells.Item(9, 81).Value = 0.04507233458906522
  $ws.Cells.Item(9, 82).Value = 0.03421069038514324
  $ws.Cells.Item(9, 83).Value = 0.02718842319100677
  $ws.Cells.Item(9, 84).Value = 0.003194083152655299
  $ws.Cells.Item(9, 85).Value = 0.004131048652010479
  $ws.Cells.Item(9, 86).Value = 0.0008464665507843835
  $ws.Cells.Item(9, 87).Value = 0.001511044199975395
  $ws.Cells.Item(9, 88).Value = 0.002719765276100597
  $ws.Cells.Item(9, 89).Value = 0.0005964622677658069
  $ws.Cells.Item(9, 90).Value = 0.002509501656698289
  $ws.Cells.Item(9, 91).Value = 0.0009051818815880179
  $ws.Cells.Item(9, 92).Value = 0.002811984509269405
  $ws.Cells.Item(9, 93).Value = 0.0007595730743912825
  $ws.Cells.Item(9, 94).Value = 0.002556805851846755
  $ws.Cells.Item(9, 95).Value = 0.003016809673591287
  $ws.Cells.Item(9, 96).Value = 0.0002974563909251212
  $ws.Cells.Item(9, 97).Value = 0.002442971905166405
  $ws.Cells.Item(9, 98).Value = 0.002167722924020633
  $ws.Cells.Item(9, 99).Value = 0.000002588713083851896
  $ws.Cells.Item(9, 100).Value = 0.002753789733325153
  $ws.Cells.Item(10, 1).Value = 0.9999235361075586
  $ws.Cells.Item(10, 2).Value = 10
  $ws.Cells.Item(10, 3).Value = 0.4015086606573693
  $ws.Cells.Item(10, 4).Value = 0.01604397074603381
  $ws.Cells.Item(10, 5).Value = 0.422281829105335
  $ws.Cells.Item(10, 6).Value = 0.3707669628383456
  $ws.Cells.Item(10, 7).Value = 0.1792748421384689
  $ws.Cells.Item(10, 8).Value = 0.006445515156503034
  $ws.Cells.Item(10, 9).Value = 0.143966692240987
  $ws.Cells.Item(10, 10).Value = 0.04729127206912419
  $ws.Cells.Item(10, 11).Value = 0.1361580941790445
  $ws.Cells.Item(10, 17).Value = 0.1843413663405821
  $ws.Cells.Item(10, 26).Value = 0.04420786344941158
  $ws.Cells.Item(10, 30).Value = 0.03709751488524476
  $ws.Cells.Item(10, 44).Value = 0.1780104094960045
  $ws.Cells.Item(10, 47).Value = 0.007312474340219226
  $ws.Cells.Item(10, 51).Value = 0.03415793323395745
  $ws.Cells.Item(10, 52).Value = 0.001736022470452692
  $ws.Cells.Item(11, 1).Value = 0.9999894140928548
  $ws.Cells.Item(11, 2).Value = 10
  $ws.Cells.Item(11, 3).Value = 0.5621231752929592
  $ws.Cells.Item(11, 4).Value = 0.01712064865790634
  $ws.Cells.Item(11, 5).Value = 0.5961625851438895
  $ws.Cells.Item(11, 6).Value = 0.5407715302800719
  $ws.Cells.Item(11, 7).Value = 0.09616565571993174
  $ws.Cells.Item(11, 8).Value = 0.004022001592029104
  $ws.Cells.Item(11, 9).Value = 0.07180985280186022
  $ws.Cells.Item(11, 10).Value = 0.02376138180011252
  $ws.Cells.Item(11, 11).Value = 0.07069688658275719
  $ws.Cells.Item(11, 17).Value = 0.103487458720009
  $ws.Cells.Item(11, 26).Value = 0.0216368817431187
  $ws.Cells.Item(11, 30).Value = 0.02070309470781175
  $ws.Cells.Item(11, 44).Value = 0.08754080281299732
  $ws.Cells.Item(11, 49).Value = 0.2231745937029853
  $ws.Cells.Item(11, 58).Value = 0.2770013898163872
